$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.903.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.80'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  -4.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3163'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07208'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08406'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7501'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.423'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.869.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.888.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.094'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("E18").Value = '  -2.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007826'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.123.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.969'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9984'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1551'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.272'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.032'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.504'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.589'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.529'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.269'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05314'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.237'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7532'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9976'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.700'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01957'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.752'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4524'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.117.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.62%  '
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8573'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.31'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.113'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.646'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.838'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.020.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.69%  '
